$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 21, shifting the existing rows (old 21-29) down to 23-31.
$ws.Range("A21:T22").Insert()

# Populate new row 21
$ws.Range("A21").Value = 2
$ws.Range("B21").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C21").Value = "Coquimbo"
$ws.Range("D21").Value = 44629
$ws.Range("E21").Value = 4
$ws.Range("F21").Value = "Fruta"
$ws.Range("G21").Value = 100103
$ws.Range("H21").Value = "Frutos de hueso (carozo)"
$ws.Range("I21").Value = 100103002
$ws.Range("J21").Value = "Ciruela"
$ws.Range("K21").Value = "Black Amber"
$ws.Range("L21").Value = "Primera"
$ws.Range("M21").Value = 20
$ws.Range("N21").Value = 220000
$ws.Range("O21").Value = 230000
$ws.Range("P21").Value = 225000
$ws.Range("Q21").Value = "$/bins (450 kilos)"
$ws.Range("R21").Value = "Región de O'Higgins"
$ws.Range("S21").Value = 500
$ws.Range("T21").Value = 450

# Populate new row 22
$ws.Range("A22").Value = 2
$ws.Range("B22").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C22").Value = "Coquimbo"
$ws.Range("D22").Value = 44629
$ws.Range("E22").Value = 4
$ws.Range("F22").Value = "Fruta"
$ws.Range("G22").Value = 100103
$ws.Range("H22").Value = "Frutos de hueso (carozo)"
$ws.Range("I22").Value = 100103002
$ws.Range("J22").Value = "Ciruela"
$ws.Range("K22").Value = "Larry Ann"
$ws.Range("L22").Value = "Primera"
$ws.Range("M22").Value = 20
$ws.Range("N22").Value = 220000
$ws.Range("O22").Value = 230000
$ws.Range("P22").Value = 225000
$ws.Range("Q22").Value = "$/bins (450 kilos)"
$ws.Range("R22").Value = "Región de O'Higgins"
$ws.Range("S22").Value = 500
$ws.Range("T22").Value = 450
